# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on the active worksheet to the newly
# recalculated strikeout (K) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 3
    21 = 2
    22 = 2
    24 = 2
    25 = 1
    26 = 3
    27 = 2
    28 = 2
    29 = 2
    30 = 3
    31 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
